# Refresh the cryptos price table (rows 2-51) with the latest scrape.
# Column D ("Price") values are numeric-looking text (e.g. "70.783.97" uses
# dots as thousands separators) so they are apostrophe-prefixed to force
# Excel to store them as text rather than silently coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = "'70.783.97"
$ws.Range('E2').Value = '  +1.00%  '

# Row 3: Ethereum
$ws.Range('D3').Value = "'3.586.80"
$ws.Range('E3').Value = '  +0.26%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  -0.05%  '

# Row 5: BNB
$ws.Range('D5').Value = "'586.51"
$ws.Range('E5').Value = '  +1.45%  '

# Row 6: Solana
$ws.Range('D6').Value = "'186.10"
$ws.Range('E6').Value = '  -0.12%  '

# Row 7: LidoStakedEther
$ws.Range('D7').Value = "'3.570.60"
$ws.Range('E7').Value = '  -0.10%  '

# Row 8: XRP
$ws.Range('D8').Value = "'0.622"
$ws.Range('E8').Value = '  +0.67%  '

# Row 9: USDC
$ws.Range('E9').Value = '  +0.11%  '

# Row 10: Dogecoin
$ws.Range('D10').Value = "'0.214"
$ws.Range('E10').Value = '  +17.19%  '

# Row 11: Cardano
$ws.Range('D11').Value = "'0.651"
$ws.Range('E11').Value = '  +0.26%  '

# Row 12: Avalanche
$ws.Range('E12').Value = '  -1.56%  '

# Row 13: ShibaInu
$ws.Range('D13').Value = "'0.0000322"
$ws.Range('E13').Value = '  +5.93%  '

# Row 14: Polkadot
$ws.Range('E14').Value = '  +0.13%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range('D15').Value = "'4.155.75"
$ws.Range('E15').Value = '  +0.10%  '

# Row 16: Chainlink
$ws.Range('D16').Value = "'19.59"
$ws.Range('E16').Value = '  -0.27%  '

# Row 17: WrappedBTC (was WrappedEther)
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = "'70.741.12"
$ws.Range('E17').Value = '  +1.02%  '

# Row 18: WrappedEther (was WrappedBTC)
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = "'3.583.85"
$ws.Range('E18').Value = '  +0.10%  '

# Row 19: Uniswap
$ws.Range('E19').Value = '  -1.23%  '

# Row 20: BitcoinCash
$ws.Range('D20').Value = "'568.54"
$ws.Range('E20').Value = '  +15.17%  '

# Row 21: TRON
$ws.Range('E21').Value = '  -0.15%  '

# Row 22: Polygon
$ws.Range('E22').Value = '  -1.79%  '

# Row 23: InternetComputer(DFINITY)
$ws.Range('D23').Value = "'17.66"
$ws.Range('E23').Value = '  -9.18%  '

# Row 24: PancakeSwap
$ws.Range('E24').Value = '  +6.12%  '

# Row 25: Toncoin
$ws.Range('E25').Value = '  -1.42%  '

# Row 26: Litecoin
$ws.Range('D26').Value = "'95.48"
$ws.Range('E26').Value = '  -1.31%  '

# Row 27: RenderToken
$ws.Range('D27').Value = "'11.48"
$ws.Range('E27').Value = '  -0.52%  '

# Row 28: ImmutableX
$ws.Range('E28').Value = '  -0.39%  '

# Row 29: Filecoin
$ws.Range('D29').Value = "'9.14"
$ws.Range('E29').Value = '  -1.88%  '

# Row 30: EthereumClassic
$ws.Range('E30').Value = '  +2.03%  '

# Row 31: NEARProtocol
$ws.Range('D31').Value = "'7.31"
$ws.Range('E31').Value = '  -5.38%  '

# Row 32: Cosmos
$ws.Range('D32').Value = "'12.41"
$ws.Range('E32').Value = '  +2.51%  '

# Row 33: OKB
$ws.Range('D33').Value = "'64.91"
$ws.Range('E33').Value = '  -1.40%  '

# Row 34: Hedera
$ws.Range('E34').Value = '  -0.84%  '

# Row 35: Bittensor (was Fetch.AI)
$ws.Range('B35').Value = 'Bittensor'
$ws.Range('C35').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D35').Value = "'563.16"
$ws.Range('E35').Value = '  -1.75%  '

# Row 36: Fetch.AI (was Bittensor)
$ws.Range('B36').Value = 'Fetch.AI'
$ws.Range('C36').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D36').Value = "'3.31"
$ws.Range('E36').Value = '  +2.19%  '

# Row 37: TheGraph
$ws.Range('D37').Value = "'0.417"
$ws.Range('E37').Value = '  +0.72%  '

# Row 38: PEPE (was InjectiveProtocol)
$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D38').Value = "'0.0₃0805"
$ws.Range('E38').Value = '  +1.83%  '

# Row 39: InjectiveProtocol (was Dai)
$ws.Range('B39').Value = 'InjectiveProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D39').Value = "'37.76"
$ws.Range('E39').Value = '  -2.63%  '

# Row 40: Dai (was PEPE)
$ws.Range('B40').Value = 'Dai'
$ws.Range('C40').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D40').Value = "'1.00"
$ws.Range('E40').Value = '  +0.07%  '

# Row 41: Maker (was dogwifhat)
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = "'3.384.81"
$ws.Range('E41').Value = '  +6.12%  '

# Row 42: dogwifhat (was Maker)
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').Value = "'3.13"
$ws.Range('E42').Value = '  -1.58%  '

# Row 43: Kaspa
$ws.Range('D43').Value = "'0.135"
$ws.Range('E43').Value = '  -0.11%  '

# Row 44: Stacks
$ws.Range('E44').Value = '  -2.46%  '

# Row 45: ApeXProtocol
$ws.Range('E45').Value = '  +0.42%  '

# Row 46: VeChain
$ws.Range('D46').Value = "'0.0446"
$ws.Range('E46').Value = '  +1.31%  '

# Row 47: ThetaToken
$ws.Range('E47').Value = '  -3.21%  '

# Row 48: THORChain
$ws.Range('E48').Value = '  -1.22%  '

# Row 49: Stellar
$ws.Range('E49').Value = '  +0.86%  '

# Row 50: FirstDigitalUSD
$ws.Range('E50').Value = '  -0.10%  '

# Row 51: OceanProtocol
$ws.Range('D51').Value = "'1.43"
$ws.Range('E51').Value = '  -8.47%  '

